# rules changes according to operation requirements
# Insert two new store rows:
#  - "受气牛肉约克郡光环店" right after "受气牛肉解放碑店" (old row 20), before the 李子坝梁山鸡 group
#  - "李子坝梁山鸡约克郡光环店" right after "李子坝梁山鸡长嘉汇店" (old row 36), before the 沸堂蛙 group

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two blank rows first (bottom one first so the first insert's
#     row index is not affected by the later one) ---

# New row for 李子坝梁山鸡约克郡光环店 goes right after old row 36 (长嘉汇店), i.e. before old row 37 (来个宝沸堂蛙全城配送店)
$ws.Rows("37:37").Insert()

# New row for 受气牛肉约克郡光环店 goes right after old row 20 (受气牛肉解放碑店), i.e. before old row 21 (李子坝梁山鸡东原悦荟店)
$ws.Rows("21:21").Insert()

# --- Fill in the values for the two new rows ---
# Row 21: 受气牛肉约克郡光环店
$ws.Range("A21").Value = "受气牛肉约克郡光环店"
$ws.Range("B21").Value = "受气牛肉约克郡光环店"
$ws.Range("C21").Value = "受气牛肉"
$ws.Range("D21").Value = 4.5
$ws.Range("E21").Value = 4.8
$ws.Range("F21").Value = 6

# Row 38: 李子坝梁山鸡约克郡光环店 (old row 37 shifted to 38 after the A21 insert above)
$ws.Range("A38").Value = "李子坝梁山鸡约克郡光环店"
$ws.Range("B38").Value = "李子坝梁山鸡约克郡光环店"
$ws.Range("C38").Value = "李子坝梁山鸡"
$ws.Range("D38").Value = 4.5
$ws.Range("E38").Value = 4.8
$ws.Range("F38").Value = 6

# --- Fix up formatting so the new "last row of group" look matches the rest of
#     the sheet: copy the style from the donor cells that already carry it ---

# A21/A38 need the "last-in-group" left-column style (already used by A2, A4, A5, A31, A39, ...)
$ws.Range("A2").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A38").PasteSpecial(-4122)

# B21/B38/C38 need style index 10 (already used by C4 before the edit)
$ws.Range("C4").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C38").PasteSpecial(-4122)

# A20 (受气牛肉解放碑店, now the second-to-last of its group) also gets the
# "last-in-group" style since the group boundary moved down by one row
$ws.Range("A2").Copy()
$ws.Range("A20").PasteSpecial(-4122)

# C37 (李子坝梁山鸡长嘉汇店, now second-to-last of the 李子坝梁山鸡 group) gets style 10 too
$ws.Range("C4").Copy()
$ws.Range("C37").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Re-set values (PasteSpecial of formats only shouldn't disturb values, but
#     make sure the brand/score values on C21/D21 etc. are still correct) ---
$ws.Range("C21").Value = "受气牛肉"
$ws.Range("D21").Value = 4.5
$ws.Range("E21").Value = 4.8
$ws.Range("F21").Value = 6

$ws.Range("C38").Value = "李子坝梁山鸡"
$ws.Range("D38").Value = 4.5
$ws.Range("E38").Value = 4.8
$ws.Range("F38").Value = 6

# --- Sheet view bookkeeping to match the target file ---
$ws.Range("I22").Select()

Write-Output "edit complete"
